$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the results table with re-run figures for the combined
# (VRM + MGN + MTS + MBL) model column, plus refreshed baseline-model
# values for Metatranscriptomics (MTS) / Metabolomics (MBL).
$ws.Range("D6").Value = "0.38 [0.15, 0.62], 0.04 "
$ws.Range("E6").Value = "0.47 [0.25, 0.69], 0.01"
$ws.Range("F2").Value = "0.80 [0.63, 0.98], 0.46"
$ws.Range("F3").Value = "3.14 [1.10, 9.00], 0.03"
$ws.Range("F5").Value = "3.13 [1.19, 8.26], 0.02"
$ws.Range("F4").Value = "0.77 [0.58, 0.96], 0.37"

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("D15").Select()
